$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: law-why
$ws.Range("A26").Value = "law-why"
$ws.Range("B26").Value = "เพื่อเป็นเครื่องมือส่งเสริมความเจริญเติบโตทางเศรษฐกิจ`n -เพื่อเป็นเครื่องมือในการควบคุมการบริโภคของประชาชนและบำรุงสาธารณูปโภคเละบริการสาธารณะ`n -เพื่อเป็นเครื่องมือในการกระจายรายได้แก่ให้ประชาชนและเป็นการรักษาเสถียรภาพในทางเศรษฐกิจของประเทศ"
$ws.Rows.Item(26).RowHeight = 102

# Row 27: law-calculate (column A has default/no explicit style)
$ws.Range("A27").Value = "law-calculate"
$ws.Range("A27").Style = "Normal"
$ws.Range("B27").Value = "การคำนวณภาษีของบุคคลธรรมดา เงินได้สุทธิซึ่งเป็นฐานภาษีสำหรับคำนวณภาษีเงินได้บุคคลธรรมดาซึ่งมาจากเงินได้พึงประเมินที่หักค่าใช้จ่ายและค่าลดหย่อนเรียบร้อยแล้ว (ค่าใช้จ่าย-ค่าลดหย่อน=เงินได้สุทธิ)"
$ws.Rows.Item(27).RowHeight = 63.75

# Row 28: law-time
$ws.Range("A28").Value = "law-time"
$ws.Range("A28").WrapText = $true
$ws.Range("B28").Value = "กฎหมายกำหนดให้บุคคลต้องทำการยื่นเสียภาษีในช่วง 1 มกราคม - 31 มีนาคม ของทุกปี"
$ws.Rows.Item(28).RowHeight = 25.5

# Row 29: law-salary (column A has default/no explicit style)
$ws.Range("A29").Value = "law-salary"
$ws.Range("A29").Style = "Normal"
$ws.Range("B29").Value = "หากมีเงินเดือนหรือมีรายได้จากหลายทางเกิน 10,000 บาท/เดือน (120,000 บาท/ปี) ต้องยื่นภาษีทุกคน"
$ws.Rows.Item(29).RowHeight = 38.25

# Update the sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C27").Select()
